$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (they already were stored as
# inline strings) so Excel does not auto-convert numeric-looking values
# like "243.18" into a float when the new value is assigned below.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Updated price / volume(1h) values
$ws.Range("D2").Value = "36.502.20"
$ws.Range("D3").Value = "1.941.03"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "243.18"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "57.07"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "0.362"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "0.0802"
$ws.Range("E10").Value = "  -2.06%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "21.74"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("D13").Value = "2.228.32"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "0.805"
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("D15").Value = "13.28"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "1.943.25"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "36.478.54"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "69.23"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0854"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "227.03"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "4.95"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -4.18%  "
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  +1.45%  "
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "159.52"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E28").Value = "  +14.66%  "
$ws.Range("D29").Value = "19.16"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").Value = "4.61"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "0.0615"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").Value = "4.17"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "6.13"
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  +3.06%  "
$ws.Range("E39").Value = "  +15.91%  "
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "0.0209"
$ws.Range("E42").Value = "  +0.44%  "
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "15.74"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "1.342.08"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "1.02"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "86.14"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "7.11"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "2.120.35"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "43.03"
$ws.Range("E51").Value = "  -6.82%  "
